$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.588.24'
$ws.Range("E2").Value = '  -2.33%  '

$ws.Range("D3").Value = '1.753.96'
$ws.Range("E3").Value = '  -3.36%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''324.44'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.43%  '

$ws.Range("D6").Value = '''0.9997'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").Value = '''0.4495'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +3.07%  '

$ws.Range("D8").Value = '''0.3617'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -1.38%  '

$ws.Range("D9").Value = '''0.07512'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -2.23%  '

$ws.Range("D10").Value = '''42.20'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -5.96%  '

$ws.Range("D11").Value = '''1.106'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -3.07%  '

$ws.Range("D12").Value = '''0.9997'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -0.11%  '

$ws.Range("D13").Value = '''20.73'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -5.82%  '

$ws.Range("D14").Value = '''6.051'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -4.19%  '

$ws.Range("D15").Value = '''7.184'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -4.18%  '

$ws.Range("D16").Value = '1.755.22'
$ws.Range("E16").Value = '  -3.54%  '

$ws.Range("D17").Value = '''92.83'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -2.35%  '

$ws.Range("E18").Value = '  -1.41%  '

$ws.Range("D19").Value = '''0.06394'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -1.38%  '

$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("D21").Value = '''16.95'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -2.58%  '

$ws.Range("D22").Value = '''5.870'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -5.91%  '

$ws.Range("D23").Value = '27.614.11'
$ws.Range("E23").Value = '  -2.34%  '

$ws.Range("D24").Value = '''11.23'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -3.06%  '

$ws.Range("D25").Value = '''2.111'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -0.99%  '

$ws.Range("D26").Value = '''161.50'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("E27").Value = '  -1.28%  '

$ws.Range("D28").Value = '1.954.03'
$ws.Range("E28").Value = '  -3.53%  '

$ws.Range("D29").Value = '''2.124'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -6.95%  '

$ws.Range("D30").Value = '''125.42'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -3.56%  '

$ws.Range("D31").Value = '''1.089'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -10.32%  '

$ws.Range("D32").Value = '''0.09043'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -1.11%  '

$ws.Range("D33").Value = '''5.573'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -7.57%  '

$ws.Range("D34").Value = '''3.636'
$ws.Range("D34").NumberFormat = "General"

$ws.Range("D35").Value = '''12.02'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -8.41%  '

$ws.Range("D36").Value = '''0.02314'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -2.84%  '

$ws.Range("D37").Value = '''0.6380'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -3.39%  '

$ws.Range("E38").Value = '  -4.14%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '''4.980'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -5.16%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '''0.05983'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -3.70%  '

$ws.Range("D41").Value = '''1.201'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.13%  '

$ws.Range("D42").Value = '''0.9999'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("D43").Value = '''1.386'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -3.00%  '

$ws.Range("D44").Value = '''7.805'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -3.28%  '

$ws.Range("D45").Value = '''13.28'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -4.10%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.5895'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -3.58%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '''3.712'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -0.77%  '

$ws.Range("D48").Value = '''1.961'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -3.14%  '

$ws.Range("D49").Value = '''121.60'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -3.36%  '

$ws.Range("D50").Value = '''1.163'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -0.11%  '

$ws.Range("D51").Value = '''0.06870'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -1.80%  '

